$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29-40 down to 30-41
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new record's data
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44900
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103003
$ws.Range("J29").Value = "Damasco"
$ws.Range("K29").Value = "Castle Brite"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 220
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 15545
$ws.Range("Q29").Value = "$/caja 10 kilos"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 1554
$ws.Range("T29").Value = 10
